$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 14: rename + new averaged-intensity values ---
$ws.Range("B14").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C14").Value = 1.000152709206348
$ws.Range("D14").Value = 1.004512157379551
$ws.Range("E14").Value = 0.9787567452273622
$ws.Range("F14").Value = 0.993029961728995
$ws.Range("G14").Value = 1.000152709206348
$ws.Range("H14").Value = 1.004512157379551
$ws.Range("I14").Value = 0.9803620705836725
$ws.Range("J14").Value = 0.981882613418525
$ws.Range("K14").Value = 0.9904940528935374
$ws.Range("L14").Value = 0.9977694582706162
$ws.Range("M14").Value = 1.000124337934872
$ws.Range("N14").Value = 0.9916344513034566
$ws.Range("O14").Value = 0.998771059554273
$ws.Range("P14").Value = 0.9944738706044204
$ws.Range("Q14").Value = 0.9920996214453028
$ws.Range("R14").Value = 0.9944738706044204
$ws.Range("S14").Value = 0.9941128933855641
$ws.Range("T14").Value = 0.995320856549721
$ws.Range("U14").Value = 0.990869971088576

# --- Append new rows 24 & 25, cloning row 23 formatting first ---
$ws.Range("A23:U23").Copy($ws.Range("A24:U24"))
$ws.Range("A23:U23").Copy($ws.Range("A25:U25"))

# Row 24
$ws.Range("A24").Value = 22
$ws.Range("B24").Value = "RotRing Axis-Y Res-5.0 Theta-2.5 "
$ws.Range("C24").Value = 1.117426683157106
$ws.Range("D24").Value = 1.005579092170347
$ws.Range("E24").Value = 0.9317258133556036
$ws.Range("F24").Value = 0.9946533588929792
$ws.Range("G24").Value = 1.117426683157106
$ws.Range("H24").Value = 1.005579092170347
$ws.Range("I24").Value = 0.9670081193522106
$ws.Range("J24").Value = 0.9429378854902224
$ws.Range("K24").Value = 1.031356459421845
$ws.Range("L24").Value = 0.985517377788908
$ws.Range("M24").Value = 1.117300236712586
$ws.Range("N24").Value = 0.9686524527629754
$ws.Range("O24").Value = 1.000116225531663
$ws.Range("P24").Value = 1.018243862894352
$ws.Range("Q24").Value = 0.9773194214729767
$ws.Range("R24").Value = 1.018243862894352
$ws.Range("S24").Value = 1.012346236894009
$ws.Range("T24").Value = 1.033362326146629
$ws.Range("U24").Value = 0.9970255987036528

# Row 25
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = "RotRing Axis-Y Res-5.0 Theta-5.0 "
$ws.Range("C25").Value = 1.099275694931267
$ws.Range("D25").Value = 0.9992939533116582
$ws.Range("E25").Value = 0.935959401678264
$ws.Range("F25").Value = 0.9908990312517337
$ws.Range("G25").Value = 1.099275694931267
$ws.Range("H25").Value = 0.9992939533116582
$ws.Range("I25").Value = 0.9720531637142232
$ws.Range("J25").Value = 0.9484581381004439
$ws.Range("K25").Value = 1.023049730476284
$ws.Range("L25").Value = 0.9832911295991779
$ws.Range("M25").Value = 1.099170278302773
$ws.Range("N25").Value = 0.9676266774949611
$ws.Range("O25").Value = 0.9950964922816959
$ws.Range("P25").Value = 1.011509683307063
$ws.Range("Q25").Value = 0.9753841287472187
$ws.Range("R25").Value = 1.011509683307063
$ws.Range("S25").Value = 1.006357020293231
$ws.Range("T25").Value = 1.024940755220838
$ws.Range("U25").Value = 0.9940350303828814
